$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF" + $row)
    if ($cell.Value2 -eq "2-25-2011-12") {
        # Leading apostrophe forces Excel to store this as literal text
        # instead of auto-parsing the ISO-like string into a date serial.
        $cell.Value2 = "'2012-02-25"
    }
}
